$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioA")

# Set ExisUnits (column F) to 0 and MaxInvest (column I) to 200 for rows 8-18
for ($r = 8; $r -le 18; $r++) {
    $ws.Cells.Item($r, 6).Value = 0
    $ws.Cells.Item($r, 9).Value = 200
}

# Update the active cell selection to match the saved view state
$ws.Activate()
$ws.Range("G23").Select()
